$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 37 - this pushes the existing rows 37..154
# down to 38..155 (matching the shift observed across the whole diff).
$ws.Rows(37).Insert()

# Populate the newly inserted row 37 with the new weekly record
# (same market/category/variety/quality template as the surrounding rows,
# with the new date and price/volume figures from the commit).
$ws.Cells.Item(37, 1).Value = 4
$ws.Cells.Item(37, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(37, 3).Value = "Los Lagos"
$ws.Cells.Item(37, 4).Value = 44481
$ws.Cells.Item(37, 5).Value = 10
$ws.Cells.Item(37, 6).Value = 100112003
$ws.Cells.Item(37, 7).Value = "Ajo"
$ws.Cells.Item(37, 8).Value = "Chino"
$ws.Cells.Item(37, 9).Value = "Primera"
$ws.Cells.Item(37, 10).Value = 400
$ws.Cells.Item(37, 11).Value = 17500
$ws.Cells.Item(37, 12).Value = 19000
$ws.Cells.Item(37, 13).Value = 18250
$ws.Cells.Item(37, 14).Value = "`$/caja 10 kilos"
$ws.Cells.Item(37, 15).Value = "China"
$ws.Cells.Item(37, 16).Value = 1825
$ws.Cells.Item(37, 17).Value = 10
$ws.Cells.Item(37, 18).Value = "Hortaliza"
